$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the header cell.
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 00:52"

# Refresh "Estados Unidos" totals (row 4).
$ws.Range("B4").Value = 611156
$ws.Range("C4").Value = 24215
$ws.Range("D4").Value = 38675
$ws.Range("E4").Value = 546557
$ws.Range("F4").Value = 13443
$ws.Range("G4").Value = 2284
$ws.Range("H4").Value = 25924

# Nigeria's case counts were refreshed, which moved it above San Marino once the
# table was re-sorted by total cases (descending). Shift San Marino / Guinea /
# Republica de Yibuti down one row and insert the updated Nigeria row above them.
$ws.Range("A107").Value = "Republica de Yibuti"
$ws.Range("B107").Value = 363
$ws.Range("C107").Value = 65
$ws.Range("D107").Value = 53
$ws.Range("E107").Value = 308
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 2

$ws.Range("A106").Value = "Guinea"
$ws.Range("B106").Value = 363
$ws.Range("C106").Value = 44
$ws.Range("D106").Value = 31
$ws.Range("E106").Value = 332
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 0

$ws.Range("A105").Value = "San Marino"
$ws.Range("B105").Value = 371
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 53
$ws.Range("E105").Value = 282
$ws.Range("F105").Value = 15
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 36

$ws.Range("A104").Value = "Nigeria"
$ws.Range("B104").Value = 373
$ws.Range("C104").Value = 30
$ws.Range("D104").Value = 99
$ws.Range("E104").Value = 263
$ws.Range("F104").Value = 2
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 11
